$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix company name typo in B3
$ws.Range("B3").Value = "New Sources Energy N.V. (ENXTAM:NSE)"

# Row 2 updates
$ws.Range("K2").Value = -0.188
$ws.Range("U2").Value = 0.043
$ws.Range("V2").Value = 0.005882352941176471
$ws.Range("W2").Value = 0.3202725724020443
$ws.Range("X2").Value = 0.04205768685451744
$ws.Range("Y2").Value = 0.2782148855475269
$ws.Range("AA2").Value = 7.782990560619196
$ws.Range("AB2").Value = 0.04203721426692397
$ws.Range("AC2").Value = 7.740953346352272
$ws.Range("AD2").Value = 0
$ws.Range("AE2").Value = 0.007654100618233418
$ws.Range("AF2").Value = 0.007654100618233418
$ws.Range("AG2").Value = -0.03534589938176658
$ws.Range("AH2").Value = 0.001045977373757904
$ws.Range("AI2").Value = 0.04592806651525017
$ws.Range("AJ2").Value = -0.004858773887099693
$ws.Range("AK2").Value = -0.2858449433140323
$ws.Range("AL2").Value = 0.005
$ws.Range("AM2").Value = 0.005
$ws.Range("AN2").Value = 0
$ws.Range("AO2").Value = -52.4
$ws.Range("AP2").Value = -8.836474845441645
$ws.Range("AQ2").Value = -52.4

# Row 3 updates
$ws.Range("K3").Value = -0.188
$ws.Range("U3").Value = 0.043
$ws.Range("V3").Value = 0.005882352941176471
$ws.Range("W3").Value = 0.3202725724020443
$ws.Range("X3").Value = 0.04205768685451744
$ws.Range("Y3").Value = 0.2782148855475269
$ws.Range("AA3").Value = 7.782990560619196
$ws.Range("AB3").Value = 0.04203721426692397
$ws.Range("AC3").Value = 7.740953346352272
$ws.Range("AD3").Value = 0
$ws.Range("AE3").Value = 0.007654100618233418
$ws.Range("AF3").Value = 0.007654100618233418
$ws.Range("AG3").Value = -0.03534589938176658
$ws.Range("AH3").Value = 0.001045977373757904
$ws.Range("AI3").Value = 0.04592806651525017
$ws.Range("AJ3").Value = -0.004858773887099693
$ws.Range("AK3").Value = -0.2858449433140323
$ws.Range("AL3").Value = 0.005
$ws.Range("AM3").Value = 0.005
$ws.Range("AN3").Value = 0
$ws.Range("AO3").Value = -52.4
$ws.Range("AP3").Value = -8.836474845441645
$ws.Range("AQ3").Value = -52.4
